$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 221 and 222: the two matches had all of their data (every column
# except the running index in column A) swapped between them.
# ---------------------------------------------------------------------------
$swapCols = @(2,6,7,8,9,10,11,12,13,14,15,16,17,18,19,21,22,23,24,26,27,28)
foreach ($c in $swapCols) {
    $v1 = $ws.Cells.Item(221, $c).Value2
    $v2 = $ws.Cells.Item(222, $c).Value2
    $ws.Cells.Item(221, $c).Value2 = $v2
    $ws.Cells.Item(222, $c).Value2 = $v1
}

# ---------------------------------------------------------------------------
# Row 239: a handful of odds columns were corrected.
# ---------------------------------------------------------------------------
$ws.Cells.Item(239, 14).Value2 = 2.25    # N239
$ws.Cells.Item(239, 16).Value2 = 3.25    # P239
$ws.Cells.Item(239, 18).Value2 = 1.975   # R239
$ws.Cells.Item(239, 19).Value2 = 1.875   # S239
$ws.Cells.Item(239, 21).Value2 = 2.025   # U239
$ws.Cells.Item(239, 22).Value2 = 1.825   # V239

# ---------------------------------------------------------------------------
# New rows 241-247 : seven newly scraped matches appended to the table.
# Columns: r, A(id), B, F(HomeTeam), G(AwayTeam), K,L,M,N,O,P,Q,R,S,T,U,V,
#          W,X,Y,Z,AA, E(date serial)
# ---------------------------------------------------------------------------
$newRows = @(
    @{ r=241; a=239; b=6803845; f='Termalica BB Nieciecza'; g='Zaglebie Sosnowiec';
       k=1.615; l=4.333; m=4.75;  n=1.615; o=4.333; p=4.75;  q=-1;    rr=2.05;  s=1.8;   t=2.75; u=1.925; v=1.925; w=0; x=0; y=0; z=0; aa=0;
       date=45388.41666666666 },
    @{ r=242; a=240; b=6803842; f='GKS Tychy 71'; g='Wisla Plock';
       k=2.25;  l=3.25; m=3.25;  n=2.25;  o=3.25; p=3.25;  q=-0.25; rr=1.95;  s=1.9;   t=2.25; u=1.8;   v=2.05;  w=0; x=0; y=0; z=0; aa=0;
       date=45388.52083333334 },
    @{ r=243; a=241; b=6803840; f='Odra Opole'; g='Miedz Legnica';
       k=2.7;   l=3.1;  m=2.7;   n=2.7;   o=3.1;  p=2.7;   q=0;     rr=1.925; s=1.925; t=2.25; u=1.95;  v=1.9;   w=0; x=0; y=0; z=0; aa=0;
       date=45388.625 },
    @{ r=244; a=242; b=6805665; f='Wisla Krakow'; g='Motor Lublin';
       k=1.571; l=4.333; m=5;    n=1.571; o=4.333; p=5;     q=-1;    rr=2.025; s=1.825; t=2.75; u=1.85;  v=2;     w=0; x=0; y=0; z=0; aa=0;
       date=45389.31944444445 },
    @{ r=245; a=243; b=6805663; f='GKS Katowice'; g='Lechia Gdansk';
       k=2.25;  l=3.5;  m=3.1;   n=2.25;  o=3.5;  p=3.1;   q=-0.25; rr=1.975; s=1.875; t=2.5;  u=2;     v=1.85;  w=0; x=0; y=0; z=0; aa=0;
       date=45389.41666666666 },
    @{ r=246; a=244; b=6803843; f='Stal Rzeszow'; g='Znicz Pruszkw';
       k=2.3;   l=2.9;  m=3.2;   n=2.375; o=3;    p=3.25;  q=-0.25; rr=2.025; s=1.825; t=2.25; u=1.925; v=1.925; w=0; x=0; y=0; z=0; aa=0;
       date=45389.54166666666 },
    @{ r=247; a=245; b=6805664; f='Polonia Warsaw'; g='Resovia Rzeszow';
       k=2.1;   l=3.5;  m=3.3;   n=2.1;   o=3.5;  p=3.3;   q=-0.25; rr=1.85;  s=2;     t=2.75; u=1.975; v=1.875; w=0; x=0; y=0; z=0; aa=0;
       date=45390.54166666666 }
)

foreach ($row in $newRows) {
    $r = $row.r

    # Column A (bold/bordered index style) - copy formatting from the row above.
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value2 = $row.a

    $ws.Cells.Item($r, 2).Value2 = $row.b
    $ws.Cells.Item($r, 3).Value2 = "Poland I Liga"
    $ws.Cells.Item($r, 4).Value2 = "Poland I Liga"

    # Column E (date/time number format) - copy formatting from the row above.
    $ws.Cells.Item($r - 1, 5).Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4122)
    $ws.Cells.Item($r, 5).Value2 = $row.date

    $ws.Cells.Item($r, 6).Value2 = $row.f
    $ws.Cells.Item($r, 7).Value2 = $row.g

    # H, I, J (FTHG/FTAG/FTR) are intentionally left blank - match not played yet.

    $ws.Cells.Item($r, 11).Value2 = $row.k
    $ws.Cells.Item($r, 12).Value2 = $row.l
    $ws.Cells.Item($r, 13).Value2 = $row.m
    $ws.Cells.Item($r, 14).Value2 = $row.n
    $ws.Cells.Item($r, 15).Value2 = $row.o
    $ws.Cells.Item($r, 16).Value2 = $row.p
    $ws.Cells.Item($r, 17).Value2 = $row.q
    $ws.Cells.Item($r, 18).Value2 = $row.rr
    $ws.Cells.Item($r, 19).Value2 = $row.s
    $ws.Cells.Item($r, 20).Value2 = $row.t
    $ws.Cells.Item($r, 21).Value2 = $row.u
    $ws.Cells.Item($r, 22).Value2 = $row.v
    $ws.Cells.Item($r, 23).Value2 = $row.w
    $ws.Cells.Item($r, 24).Value2 = $row.x
    $ws.Cells.Item($r, 25).Value2 = $row.y
    $ws.Cells.Item($r, 26).Value2 = $row.z
    $ws.Cells.Item($r, 27).Value2 = $row.aa
}
